$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '66.676.75'
$ws.Range("E2").Value = '  +0.64%  '
$ws.Range("D3").Value = '3.619.65'
$ws.Range("E3").Value = '  +1.42%  '
$ws.Range("E4").Value = '  -0.08%  '
$cell = $ws.Range("D5")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '610.35'
$cell.Style = $origStyle
$ws.Range("E5").Value = '  +0.27%  '
$cell = $ws.Range("D6")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '150.75'
$cell.Style = $origStyle
$ws.Range("E6").Value = '  +3.67%  '
$ws.Range("D7").Value = '3.617.21'
$ws.Range("E7").Value = '  +1.42%  '
$ws.Range("E8").Value = '  -0.05%  '
$ws.Range("E9").Value = '  -0.58%  '
$ws.Range("E10").Value = '  +0.14%  '
$ws.Range("E11").Value = '  +0.35%  '
$ws.Range("E12").Value = '  +0.56%  '
$ws.Range("D13").Value = '4.232.18'
$ws.Range("E13").Value = '  +1.36%  '
$ws.Range("E14").Value = '  +0.66%  '
$cell = $ws.Range("D15")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '29.99'
$cell.Style = $origStyle
$ws.Range("E15").Value = '  +0.04%  '
$ws.Range("D16").Value = '3.620.51'
$ws.Range("E16").Value = '  +1.34%  '
$ws.Range("D17").Value = '66.746.61'
$ws.Range("E17").Value = '  +0.59%  '
$ws.Range("E18").Value = '  +1.70%  '
$ws.Range("E19").Value = '  +1.52%  '
$cell = $ws.Range("D20")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '6.38'
$cell.Style = $origStyle
$ws.Range("E20").Value = '  +2.40%  '
$cell = $ws.Range("D21")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '15.12'
$cell.Style = $origStyle
$ws.Range("E21").Value = '  +1.64%  '
$cell = $ws.Range("D22")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '428.07'
$cell.Style = $origStyle
$ws.Range("E22").Value = '  -0.41%  '
$ws.Range("E24").Value = '  -0.44%  '
$ws.Range("E25").Value = '  -0.02%  '
$ws.Range("E26").Value = '  +4.16%  '
$cell = $ws.Range("D27")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '8.38'
$cell.Style = $origStyle
$ws.Range("E27").Value = '  +5.52%  '
$ws.Range("E28").Value = '  +5.76%  '
$ws.Range("E29").Value = '  +0.44%  '
$ws.Range("E30").Value = '  +0.09%  '
$ws.Range("D31").Value = '3.616.42'
$ws.Range("E31").Value = '  +1.41%  '
$cell = $ws.Range("D32")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.159'
$cell.Style = $origStyle
$ws.Range("E32").Value = '  +3.51%  '
$ws.Range("E33").Value = '  +0.92%  '
$ws.Range("E34").Value = '  -0.71%  '
$ws.Range("E35").Value = '  +0.26%  '
$ws.Range("E36").Value = '  +0.00%  '
$cell = $ws.Range("D37")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '5.63'
$cell.Style = $origStyle
$ws.Range("E37").Value = '  +0.47%  '
$ws.Range("E38").Value = '  -1.79%  '
$cell = $ws.Range("D39")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '177.22'
$cell.Style = $origStyle
$ws.Range("E39").Value = '  -0.08%  '
$cell = $ws.Range("D40")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.0863'
$cell.Style = $origStyle
$ws.Range("E40").Value = '  +1.79%  '
$ws.Range("E41").Value = '  +0.16%  '
$cell = $ws.Range("D42")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.900'
$cell.Style = $origStyle
$ws.Range("E42").Value = '  +0.39%  '
$ws.Range("E43").Value = '  -2.27%  '
$cell = $ws.Range("D44")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '46.27'
$cell.Style = $origStyle
$ws.Range("E44").Value = '  +0.19%  '
$cell = $ws.Range("D45")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '2.60'
$cell.Style = $origStyle
$ws.Range("E45").Value = '  +7.42%  '
$cell = $ws.Range("D46")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = $origStyle
$ws.Range("E46").Value = '  -0.06%  '
$ws.Range("E47").Value = '  -2.45%  '
$ws.Range("E48").Value = '  -3.46%  '
$ws.Range("E49").Value = '  +1.45%  '
$ws.Range("E50").Value = '  +0.80%  '
$cell = $ws.Range("D51")
$origStyle = $cell.Style
$cell.NumberFormat = "@"
$cell.Value = '0.967'
$cell.Style = $origStyle
$ws.Range("E51").Value = '  +2.50%  '
